$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'71.860.79"
$ws.Range('E2').Value = "  +3.36%  "
$ws.Range('D3').Value = "'3.632.88"
$ws.Range('E3').Value = "  +7.11%  "
$ws.Range('E4').Value = "  +0.10%  "
$ws.Range('D5').Value = "'587.73"
$ws.Range('E5').Value = "  +0.46%  "
$ws.Range('D6').Value = "'181.60"
$ws.Range('E6').Value = "  +0.76%  "
$ws.Range('D7').Value = "'3.629.33"
$ws.Range('E7').Value = "  +7.20%  "
$ws.Range('D8').Value = "'0.617"
$ws.Range('E8').Value = "  +3.36%  "
$ws.Range('D9').Value = "'0.999"
$ws.Range('E9').Value = "  -0.07%  "
$ws.Range('E10').Value = "  +1.21%  "
$ws.Range('E11').Value = "  +2.82%  "
$ws.Range('D12').Value = "'49.71"
$ws.Range('E12').Value = "  +2.87%  "
$ws.Range('E13').Value = "  +0.23%  "
$ws.Range('D14').Value = "'683.89"
$ws.Range('E14').Value = "  +0.00%  "
$ws.Range('D15').Value = "'4.217.53"
$ws.Range('E15').Value = "  +6.91%  "
$ws.Range('E16').Value = "  +4.25%  "
$ws.Range('B17').Value = "WrappedEther"
$ws.Range('C17').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D17').Value = "'3.656.17"
$ws.Range('E17').Value = "  +7.81%  "
$ws.Range('B18').Value = "WrappedBTC"
$ws.Range('C18').Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('D18').Value = "'71.881.19"
$ws.Range('E18').Value = "  +3.41%  "
$ws.Range('E19').Value = "  +1.75%  "
$ws.Range('D20').Value = "'18.30"
$ws.Range('E20').Value = "  +3.34%  "
$ws.Range('E21').Value = "  +2.47%  "
$ws.Range('D22').Value = "'0.941"
$ws.Range('E22').Value = "  +3.23%  "
$ws.Range('D23').Value = "'5.93"
$ws.Range('E23').Value = "  +10.90%  "
$ws.Range('D24').Value = "'17.80"
$ws.Range('E24').Value = "  +3.14%  "
$ws.Range('D25').Value = "'103.20"
$ws.Range('E25').Value = "  +0.63%  "
$ws.Range('D26').Value = "'4.01"
$ws.Range('E26').Value = "  +2.03%  "
$ws.Range('D27').Value = "'2.84"
$ws.Range('E27').Value = "  +4.84%  "
$ws.Range('E28').Value = "  +3.14%  "
$ws.Range('D29').Value = "'35.24"
$ws.Range('E29').Value = "  +3.91%  "
$ws.Range('E30').Value = "  +4.72%  "
$ws.Range('D31').Value = "'7.35"
$ws.Range('E31').Value = "  +5.89%  "
$ws.Range('D32').Value = "'4.24"
$ws.Range('E32').Value = "  +16.57%  "
$ws.Range('D33').Value = "'584.92"
$ws.Range('E33').Value = "  +5.66%  "
$ws.Range('E34').Value = "  +1.83%  "
$ws.Range('E35').Value = "  +1.80%  "
$ws.Range('D36').Value = "'59.39"
$ws.Range('E36').Value = "  +1.81%  "
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = "  -0.04%  "
$ws.Range('E38').Value = "  +0.25%  "
$ws.Range('E39').Value = "  +2.01%  "
$ws.Range('D40').Value = "'35.73"
$ws.Range('E40').Value = "  +0.11%  "
$ws.Range('D41').Value = "'0.0₃0763"
$ws.Range('E41').Value = "  +4.63%  "
$ws.Range('E42').Value = "  +4.04%  "
$ws.Range('E43').Value = "  +9.20%  "
$ws.Range('E44').Value = "  +2.79%  "
$ws.Range('D45').Value = "'0.346"
$ws.Range('E45').Value = "  +2.18%  "
$ws.Range('E46').Value = "  +2.46%  "
$ws.Range('D47').Value = "'2.82"
$ws.Range('E47').Value = "  +5.53%  "
$ws.Range('D48').Value = "'0.133"
$ws.Range('E48').Value = "  +2.98%  "
$ws.Range('E49').Value = "  +3.81%  "
$ws.Range('D50').Value = "'0.999"
$ws.Range('E50').Value = "  -0.07%  "
$ws.Range('D51').Value = "'131.60"
$ws.Range('E51').Value = "  +1.45%  "
